$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Baza podataka" cell to "Analiza sadržaja"
$ws.Range("B3").Value = "Analiza sadržaja"

# Update selection to match the new active cell/selection
$ws.Range("B3").Select()
